$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 144588
$ws.Range("E2").Value = 5288
$ws.Range("F2").Value = 5288
$ws.Range("G2").Value = 5309
$ws.Range("H2").Value = 4293
$ws.Range("I2").Value = 4221
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 349321
$ws.Range("L2").Value = 313409
$ws.Range("M2").Value = 35912
$ws.Range("N2").Value = 34537
$ws.Range("O2").Value = 1375
$ws.Range("P2").Value = 354
$ws.Range("Q2").Value = 27405
$ws.Range("R2").Value = -23885
$ws.Range("S2").Value = -1211
$ws.Range("T2").Value = 715
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 3.66
$ws.Range("X2").Value = 2.97
$ws.Range("Y2").Value = 13.92
$ws.Range("Z2").Value = 1.32
$ws.Range("AA2").Value = 872.72
$ws.Range("AB2").Value = 10128.34
$ws.Range("AC2").Value = 5961
$ws.Range("AD2").Value = 9.23
$ws.Range("AE2").Value = 54562
$ws.Range("AF2").Value = 1.01
$ws.Range("AG2").Value = 1450
$ws.Range("AH2").Value = 2.64
$ws.Range("AI2").Value = 21.75
$ws.Range("AJ2").Value = 70800000
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 159562
$ws.Range("E3").Value = 5662
$ws.Range("F3").Value = 5662
$ws.Range("G3").Value = 5455
$ws.Range("H3").Value = 4304
$ws.Range("I3").Value = 4188
$ws.Range("J3").Value = 116
$ws.Range("K3").Value = 400184
$ws.Range("L3").Value = 360514
$ws.Range("M3").Value = 39669
$ws.Range("N3").Value = 38058
$ws.Range("O3").Value = 1611
$ws.Range("P3").Value = 354
$ws.Range("Q3").Value = 24331
$ws.Range("R3").Value = -20659
$ws.Range("S3").Value = -812
$ws.Range("T3").Value = 633
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 3.55
$ws.Range("X3").Value = 2.7
$ws.Range("Y3").Value = 11.86
$ws.Range("Z3").Value = 1.15
$ws.Range("AA3").Value = 908.8
$ws.Range("AB3").Value = 11189.83
$ws.Range("AC3").Value = 5915
$ws.Range("AD3").Value = 11.88
$ws.Range("AE3").Value = 60125
$ws.Range("AF3").Value = 1.17
$ws.Range("AG3").Value = 1550
$ws.Range("AH3").Value = 2.2
$ws.Range("AI3").Value = 23.43
$ws.Range("AJ3").Value = 70800000
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 170672
$ws.Range("E4").Value = 7261
$ws.Range("F4").Value = 7261
$ws.Range("G4").Value = 7037
$ws.Range("H4").Value = 5338
$ws.Range("I4").Value = 5237
$ws.Range("J4").Value = 101
$ws.Range("K4").Value = 448634
$ws.Range("L4").Value = 405187
$ws.Range("M4").Value = 43447
$ws.Range("N4").Value = 41697
$ws.Range("O4").Value = 1751
$ws.Range("P4").Value = 354
$ws.Range("Q4").Value = 31026
$ws.Range("R4").Value = -30135
$ws.Range("S4").Value = -940
$ws.Range("T4").Value = 218
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 4.25
$ws.Range("X4").Value = 3.13
$ws.Range("Y4").Value = 13.39
$ws.Range("Z4").Value = 1.26
$ws.Range("AA4").Value = 932.59
$ws.Range("AB4").Value = 12257.02
$ws.Range("AC4").Value = 7396
$ws.Range("AD4").Value = 8.449999999999999
$ws.Range("AE4").Value = 65873
$ws.Range("AF4").Value = 0.95
$ws.Range("AG4").Value = 1650
$ws.Range("AH4").Value = 2.64
$ws.Range("AI4").Value = 19.94
$ws.Range("AJ4").Value = 70800000
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 178553
$ws.Range("E5").Value = 8679
$ws.Range("F5").Value = 8679
$ws.Range("G5").Value = 8893
$ws.Range("H5").Value = 6692
$ws.Range("I5").Value = 6611
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 483756
$ws.Range("L5").Value = 435426
$ws.Range("M5").Value = 48330
$ws.Range("N5").Value = 46280
$ws.Range("O5").Value = 2050
$ws.Range("P5").Value = 354
$ws.Range("Q5").Value = -3160
$ws.Range("R5").Value = -6042
$ws.Range("S5").Value = 5011
$ws.Range("T5").Value = 203
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 4.86
$ws.Range("X5").Value = 3.75
$ws.Range("Y5").Value = 15.21
$ws.Range("Z5").Value = 1.44
$ws.Range("AA5").Value = 900.9400000000001
$ws.Range("AB5").Value = 13636.34
$ws.Range("AC5").Value = 9338
$ws.Range("AD5").Value = 7.62
$ws.Range("AE5").Value = 73115
$ws.Range("AF5").Value = 0.97
$ws.Range("AG5").Value = 2300
$ws.Range("AH5").Value = 3.23
$ws.Range("AI5").Value = 22.02
$ws.Range("AJ5").Value = 70800000
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 174545
$ws.Range("E6").Value = 7207
$ws.Range("F6").Value = 7207
$ws.Range("G6").Value = 7353
$ws.Range("H6").Value = 5378
$ws.Range("I6").Value = 5325
$ws.Range("K6").Value = 509291
$ws.Range("L6").Value = 456363
$ws.Range("M6").Value = 52928
$ws.Range("N6").Value = 50978
$ws.Range("P6").Value = 354
$ws.Range("Q6").Value = 2464
$ws.Range("R6").Value = 1445
$ws.Range("S6").Value = -448
$ws.Range("T6").Value = 226
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 4.13
$ws.Range("X6").Value = 3.08
$ws.Range("Y6").Value = 11.06
$ws.Range("Z6").Value = 1.08
$ws.Range("AA6").Value = 862.24
$ws.Range("AB6").Value = 14935.07
$ws.Range("AC6").Value = 7522
$ws.Range("AD6").Value = 9.359999999999999
$ws.Range("AE6").Value = 80536
$ws.Range("AF6").Value = 0.87
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 2.84
$ws.Range("AI6").Value = 23.77
$ws.Range("AJ6").Value = 70800000
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 118980
$ws.Range("E7").Value = 5460
$ws.Range("G7").Value = 5350
$ws.Range("H7").Value = 3926
$ws.Range("I7").Value = 3770
$ws.Range("K7").Value = 497201
$ws.Range("L7").Value = 434162
$ws.Range("M7").Value = 63044
$ws.Range("N7").Value = 64805
$ws.Range("P7").Value = 352
$ws.Range("W7").Value = 4.59
$ws.Range("X7").Value = 3.3
$ws.Range("Y7").Value = 6.51
$ws.Range("Z7").Value = 0.78
$ws.Range("AA7").Value = 688.67
$ws.Range("AC7").Value = 5325
$ws.Range("AD7").Value = 8.359999999999999
$ws.Range("AE7").Value = 102380
$ws.Range("AF7").Value = 0.43
$ws.Range("AG7").Value = 1689
$ws.Range("AH7").Value = 3.79
$ws.Range("AI7").Value = 31.71
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 124596
$ws.Range("E8").Value = 6110
$ws.Range("G8").Value = 6050
$ws.Range("H8").Value = 4452
$ws.Range("I8").Value = 4540
$ws.Range("K8").Value = 526940
$ws.Range("L8").Value = 460766
$ws.Range("M8").Value = 66172
$ws.Range("N8").Value = 67955
$ws.Range("P8").Value = 352
$ws.Range("W8").Value = 4.9
$ws.Range("X8").Value = 3.57
$ws.Range("Y8").Value = 6.84
$ws.Range("Z8").Value = 0.87
$ws.Range("AA8").Value = 696.3200000000001
$ws.Range("AC8").Value = 6412
$ws.Range("AD8").Value = 6.63
$ws.Range("AE8").Value = 107357
$ws.Range("AF8").Value = 0.4
$ws.Range("AG8").Value = 1786
$ws.Range("AH8").Value = 4.2
$ws.Range("AI8").Value = 27.85
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 128452
$ws.Range("E9").Value = 7070
$ws.Range("G9").Value = 6982
$ws.Range("H9").Value = 5132
$ws.Range("I9").Value = 5100
$ws.Range("K9").Value = 556639
$ws.Range("L9").Value = 486783
$ws.Range("M9").Value = 69856
$ws.Range("N9").Value = 71535
$ws.Range("P9").Value = 352
$ws.Range("W9").Value = 5.5
$ws.Range("X9").Value = 4
$ws.Range("Y9").Value = 7.31
$ws.Range("Z9").Value = 0.95
$ws.Range("AA9").Value = 696.84
$ws.Range("AC9").Value = 7203
$ws.Range("AD9").Value = 5.9
$ws.Range("AE9").Value = 113012
$ws.Range("AF9").Value = 0.38
$ws.Range("AG9").Value = 2047
$ws.Range("AH9").Value = 4.82
$ws.Range("AI9").Value = 28.42
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
